$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture original values for the columns that change (D, L, M, N, O, P, Q, S, T)
# before writing, since rows 2-20 are being permuted amongst each other.
$d2 = 44966
$l2 = "Primera"
$m2 = 4
$n2 = 250000
$o2 = 250000
$p2 = 250000
$q2 = "`$/bins (400 kilos)"
$s2 = 625
$t2 = 400

$d3 = 44966
$l3 = "Primera"
$m3 = 80
$n3 = 15000
$o3 = 15000
$p3 = 15000
$q3 = "`$/caja 18 kilos granel"
$s3 = 833
$t3 = 18

$d4 = 45079
$l4 = "Primera"
$m4 = 100
$n4 = 18000
$o4 = 18000
$p4 = 18000
$q4 = "`$/caja 18 kilos granel"
$s4 = 1000
$t4 = 18

$d5 = 44411
$l5 = "Primera"
$m5 = 210
$n5 = 8000
$o5 = 8000
$p5 = 8000
$q5 = "`$/bandeja 8 kilos"
$s5 = 1000
$t5 = 8

$d6 = 45085
$l6 = "Primera"
$m6 = 110
$n6 = 16000
$o6 = 16000
$p6 = 16000
$q6 = "`$/caja 18 kilos granel"
$s6 = 889
$t6 = 18

$d7 = 44208
$l7 = "Especial"
$m7 = 70
$n7 = 24000
$o7 = 24000
$p7 = 24000
$q7 = "`$/caja 15 kilos granel"
$s7 = 1600
$t7 = 15

$d8 = 45083
$l8 = "Primera"
$m8 = 55
$n8 = 16000
$o8 = 16000
$p8 = 16000
$q8 = "`$/caja 18 kilos granel"
$s8 = 889
$t8 = 18

$d9 = 44601
$l9 = "Primera"
$m9 = 30
$n9 = 28000
$o9 = 28000
$p9 = 28000
$q9 = "`$/caja 18 kilos granel"
$s9 = 1556
$t9 = 18

$d10 = 45090
$l10 = "Primera"
$m10 = 140
$n10 = 16000
$o10 = 16000
$p10 = 16000
$q10 = "`$/caja 18 kilos granel"
$s10 = 889
$t10 = 18

$d11 = 44264
$l11 = "Calibre 100"
$m11 = 50
$n11 = 20000
$o11 = 20000
$p11 = 20000
$q11 = "`$/caja 18 kilos embalada"
$s11 = 1111
$t11 = 18

$d12 = 45086
$l12 = "Primera"
$m12 = 80
$n12 = 16000
$o12 = 16000
$p12 = 16000
$q12 = "`$/caja 18 kilos granel"
$s12 = 889
$t12 = 18

$d13 = 45093
$l13 = "Primera"
$m13 = 170
$n13 = 15000
$o13 = 16000
$p13 = 15471
$q13 = "`$/caja 18 kilos granel"
$s13 = 860
$t13 = 18

$d14 = 44511
$l14 = "Primera"
$m14 = 15
$n14 = 22000
$o14 = 22000
$p14 = 22000
$q14 = "`$/caja 15 kilos granel"
$s14 = 1467
$t14 = 15

$d15 = 45092
$l15 = "Primera"
$m15 = 220
$n15 = 16000
$o15 = 16000
$p15 = 16000
$q15 = "`$/caja 18 kilos granel"
$s15 = 889
$t15 = 18

$d16 = 44392
$l16 = "Especial"
$m16 = 500
$n16 = 7000
$o16 = 7000
$p16 = 7000
$q16 = "`$/bandeja 8 kilos"
$s16 = 875
$t16 = 8

$d17 = 44217
$l17 = "Primera"
$m17 = 55
$n17 = 18000
$o17 = 18000
$p17 = 18000
$q17 = "`$/caja 18 kilos granel"
$s17 = 1000
$t17 = 18

$d18 = 45089
$l18 = "Primera"
$m18 = 100
$n18 = 16000
$o18 = 16000
$p18 = 16000
$q18 = "`$/caja 18 kilos granel"
$s18 = 889
$t18 = 18

$d19 = 44418
$l19 = "Especial"
$m19 = 100
$n19 = 8000
$o19 = 8000
$p19 = 8000
$q19 = "`$/caja 15 kilos granel"
$s19 = 533
$t19 = 15

$d20 = 44427
$l20 = "Primera"
$m20 = 55
$n20 = 7000
$o20 = 7000
$p20 = 7000
$q20 = "`$/caja 15 kilos granel"
$s20 = 467
$t20 = 15

# Apply the new (post-edit) values: each destination row receives the
# captured values of the corresponding source row from the permutation.

$ws.Range("D2").Value = $d20
$ws.Range("L2").Value = $l20
$ws.Range("M2").Value = $m20
$ws.Range("N2").Value = $n20
$ws.Range("O2").Value = $o20
$ws.Range("P2").Value = $p20
$ws.Range("Q2").Value = $q20
$ws.Range("S2").Value = $s20
$ws.Range("T2").Value = $t20

$ws.Range("D3").Value = $d6
$ws.Range("L3").Value = $l6
$ws.Range("M3").Value = $m6
$ws.Range("N3").Value = $n6
$ws.Range("O3").Value = $o6
$ws.Range("P3").Value = $p6
$ws.Range("Q3").Value = $q6
$ws.Range("S3").Value = $s6
$ws.Range("T3").Value = $t6

$ws.Range("D4").Value = $d19
$ws.Range("L4").Value = $l19
$ws.Range("M4").Value = $m19
$ws.Range("N4").Value = $n19
$ws.Range("O4").Value = $o19
$ws.Range("P4").Value = $p19
$ws.Range("Q4").Value = $q19
$ws.Range("S4").Value = $s19
$ws.Range("T4").Value = $t19

$ws.Range("D5").Value = $d15
$ws.Range("L5").Value = $l15
$ws.Range("M5").Value = $m15
$ws.Range("N5").Value = $n15
$ws.Range("O5").Value = $o15
$ws.Range("P5").Value = $p15
$ws.Range("Q5").Value = $q15
$ws.Range("S5").Value = $s15
$ws.Range("T5").Value = $t15

$ws.Range("D6").Value = $d7
$ws.Range("L6").Value = $l7
$ws.Range("M6").Value = $m7
$ws.Range("N6").Value = $n7
$ws.Range("O6").Value = $o7
$ws.Range("P6").Value = $p7
$ws.Range("Q6").Value = $q7
$ws.Range("S6").Value = $s7
$ws.Range("T6").Value = $t7

$ws.Range("D7").Value = $d5
$ws.Range("L7").Value = $l5
$ws.Range("M7").Value = $m5
$ws.Range("N7").Value = $n5
$ws.Range("O7").Value = $o5
$ws.Range("P7").Value = $p5
$ws.Range("Q7").Value = $q5
$ws.Range("S7").Value = $s5
$ws.Range("T7").Value = $t5

$ws.Range("D8").Value = $d10
$ws.Range("L8").Value = $l10
$ws.Range("M8").Value = $m10
$ws.Range("N8").Value = $n10
$ws.Range("O8").Value = $o10
$ws.Range("P8").Value = $p10
$ws.Range("Q8").Value = $q10
$ws.Range("S8").Value = $s10
$ws.Range("T8").Value = $t10

$ws.Range("D9").Value = $d4
$ws.Range("L9").Value = $l4
$ws.Range("M9").Value = $m4
$ws.Range("N9").Value = $n4
$ws.Range("O9").Value = $o4
$ws.Range("P9").Value = $p4
$ws.Range("Q9").Value = $q4
$ws.Range("S9").Value = $s4
$ws.Range("T9").Value = $t4

$ws.Range("D10").Value = $d18
$ws.Range("L10").Value = $l18
$ws.Range("M10").Value = $m18
$ws.Range("N10").Value = $n18
$ws.Range("O10").Value = $o18
$ws.Range("P10").Value = $p18
$ws.Range("Q10").Value = $q18
$ws.Range("S10").Value = $s18
$ws.Range("T10").Value = $t18

$ws.Range("D11").Value = $d16
$ws.Range("L11").Value = $l16
$ws.Range("M11").Value = $m16
$ws.Range("N11").Value = $n16
$ws.Range("O11").Value = $o16
$ws.Range("P11").Value = $p16
$ws.Range("Q11").Value = $q16
$ws.Range("S11").Value = $s16
$ws.Range("T11").Value = $t16

$ws.Range("D12").Value = $d9
$ws.Range("L12").Value = $l9
$ws.Range("M12").Value = $m9
$ws.Range("N12").Value = $n9
$ws.Range("O12").Value = $o9
$ws.Range("P12").Value = $p9
$ws.Range("Q12").Value = $q9
$ws.Range("S12").Value = $s9
$ws.Range("T12").Value = $t9

$ws.Range("D13").Value = $d17
$ws.Range("L13").Value = $l17
$ws.Range("M13").Value = $m17
$ws.Range("N13").Value = $n17
$ws.Range("O13").Value = $o17
$ws.Range("P13").Value = $p17
$ws.Range("Q13").Value = $q17
$ws.Range("S13").Value = $s17
$ws.Range("T13").Value = $t17

$ws.Range("D14").Value = $d8
$ws.Range("L14").Value = $l8
$ws.Range("M14").Value = $m8
$ws.Range("N14").Value = $n8
$ws.Range("O14").Value = $o8
$ws.Range("P14").Value = $p8
$ws.Range("Q14").Value = $q8
$ws.Range("S14").Value = $s8
$ws.Range("T14").Value = $t8

$ws.Range("D15").Value = $d2
$ws.Range("L15").Value = $l2
$ws.Range("M15").Value = $m2
$ws.Range("N15").Value = $n2
$ws.Range("O15").Value = $o2
$ws.Range("P15").Value = $p2
$ws.Range("Q15").Value = $q2
$ws.Range("S15").Value = $s2
$ws.Range("T15").Value = $t2

$ws.Range("D16").Value = $d3
$ws.Range("L16").Value = $l3
$ws.Range("M16").Value = $m3
$ws.Range("N16").Value = $n3
$ws.Range("O16").Value = $o3
$ws.Range("P16").Value = $p3
$ws.Range("Q16").Value = $q3
$ws.Range("S16").Value = $s3
$ws.Range("T16").Value = $t3

$ws.Range("D17").Value = $d11
$ws.Range("L17").Value = $l11
$ws.Range("M17").Value = $m11
$ws.Range("N17").Value = $n11
$ws.Range("O17").Value = $o11
$ws.Range("P17").Value = $p11
$ws.Range("Q17").Value = $q11
$ws.Range("S17").Value = $s11
$ws.Range("T17").Value = $t11

$ws.Range("D18").Value = $d12
$ws.Range("L18").Value = $l12
$ws.Range("M18").Value = $m12
$ws.Range("N18").Value = $n12
$ws.Range("O18").Value = $o12
$ws.Range("P18").Value = $p12
$ws.Range("Q18").Value = $q12
$ws.Range("S18").Value = $s12
$ws.Range("T18").Value = $t12

$ws.Range("D19").Value = $d14
$ws.Range("L19").Value = $l14
$ws.Range("M19").Value = $m14
$ws.Range("N19").Value = $n14
$ws.Range("O19").Value = $o14
$ws.Range("P19").Value = $p14
$ws.Range("Q19").Value = $q14
$ws.Range("S19").Value = $s14
$ws.Range("T19").Value = $t14

$ws.Range("D20").Value = $d13
$ws.Range("L20").Value = $l13
$ws.Range("M20").Value = $m13
$ws.Range("N20").Value = $n13
$ws.Range("O20").Value = $o13
$ws.Range("P20").Value = $p13
$ws.Range("Q20").Value = $q13
$ws.Range("S20").Value = $s13
$ws.Range("T20").Value = $t13

